$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy style from existing header cell (AC1) to preserve bold/border/alignment formatting
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-40: Yankees 1998 season record (114 wins, 48 losses, 0 ties) repeated for every player row
for ($row = 2; $row -le 40; $row++) {
    $ws.Cells.Item($row, 30).Value = 114
    $ws.Cells.Item($row, 31).Value = 48
    $ws.Cells.Item($row, 32).Value = 0
}

Write-Host "done"
